$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.954938411712646
$ws.Range("B1").Value = 1.922278046607971
$ws.Range("C1").Value = 1.866639494895935
$ws.Range("D1").Value = 2.765007972717285
$ws.Range("E1").Value = 4.965007305145264
